# Actualización 11 de Mayo - Mañana
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Estadisticos 1P"
# Columns: A=Mat B=Grupo C=Totales D=Blancos E=Reprobados F=Aprobados
#          G=Por_Apro H=Promedio
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Cells.Item(2, 4).Value = 9
$ws1.Cells.Item(2, 6).Value = 20
$ws1.Cells.Item(2, 7).Value = 68.97
$ws1.Cells.Item(2, 8).Value = 8.6

$ws1.Cells.Item(3, 4).Value = 8
$ws1.Cells.Item(3, 6).Value = 21
$ws1.Cells.Item(3, 7).Value = 72.41
$ws1.Cells.Item(3, 8).Value = 8.6

$ws1.Cells.Item(4, 4).Value = 4
$ws1.Cells.Item(4, 6).Value = 18
$ws1.Cells.Item(4, 7).Value = 81.81999999999999
$ws1.Cells.Item(4, 8).Value = 8.300000000000001

# ---------------------------------------------------------------------
# Sheet 2: "Estadisticos 2P"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Cells.Item(2, 4).Value = 9
$ws2.Cells.Item(2, 5).Value = 0
$ws2.Cells.Item(2, 6).Value = 20
$ws2.Cells.Item(2, 7).Value = 68.97
$ws2.Cells.Item(2, 8).Value = 8.6

$ws2.Cells.Item(3, 4).Value = 10
$ws2.Cells.Item(3, 5).Value = 2
$ws2.Cells.Item(3, 6).Value = 19
$ws2.Cells.Item(3, 7).Value = 65.52
$ws2.Cells.Item(3, 8).Value = 8.699999999999999

$ws2.Cells.Item(4, 4).Value = 5
$ws2.Cells.Item(4, 5).Value = 1
$ws2.Cells.Item(4, 6).Value = 17
$ws2.Cells.Item(4, 7).Value = 77.27
$ws2.Cells.Item(4, 8).Value = 8.4

# ---------------------------------------------------------------------
# Sheet 3: "Estadisticos Final"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Cells.Item(2, 4).Value = 9
$ws3.Cells.Item(2, 6).Value = 20
$ws3.Cells.Item(2, 7).Value = 68.97
$ws3.Cells.Item(2, 8).Value = 8.699999999999999

$ws3.Cells.Item(3, 4).Value = 8
$ws3.Cells.Item(3, 6).Value = 21
$ws3.Cells.Item(3, 7).Value = 72.41
$ws3.Cells.Item(3, 8).Value = 8.300000000000001

$ws3.Cells.Item(4, 4).Value = 4
$ws3.Cells.Item(4, 6).Value = 18
$ws3.Cells.Item(4, 7).Value = 81.81999999999999
$ws3.Cells.Item(4, 8).Value = 8

# ---------------------------------------------------------------------
# Sheet 4: "Rescatables" -- remove the rescued-students rows (2-4),
# keeping only the header row.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Range("A2:G4").Delete()
